$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 13, pushing the existing rows 13-32 down to 15-34.
$ws.Rows("13:14").Insert()

# New row 13: Packham's Triumph, Segunda, $/caja 18 kilos granel
$ws.Cells.Item(13,1).Value = 1
$ws.Cells.Item(13,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(13,3).Value = "Arica y Parinacota"
$ws.Cells.Item(13,4).Value = 44763
$ws.Cells.Item(13,5).Value = 15
$ws.Cells.Item(13,6).Value = "Fruta"
$ws.Cells.Item(13,7).Value = 100104
$ws.Cells.Item(13,8).Value = "Frutos de pepita"
$ws.Cells.Item(13,9).Value = 100104005
$ws.Cells.Item(13,10).Value = "Pera"
$ws.Cells.Item(13,11).Value = "Packham's Triumph"
$ws.Cells.Item(13,12).Value = "Segunda"
$ws.Cells.Item(13,13).Value = 300
$ws.Cells.Item(13,14).Value = 17000
$ws.Cells.Item(13,15).Value = 18000
$ws.Cells.Item(13,16).Value = 17500
$ws.Cells.Item(13,17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(13,18).Value = "Región de O'Higgins"
$ws.Cells.Item(13,19).Value = 972
$ws.Cells.Item(13,20).Value = 18

# New row 14: Winter Nelis, Segunda, $/bandeja 18 kilos granel
$ws.Cells.Item(14,1).Value = 1
$ws.Cells.Item(14,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(14,3).Value = "Arica y Parinacota"
$ws.Cells.Item(14,4).Value = 44763
$ws.Cells.Item(14,5).Value = 15
$ws.Cells.Item(14,6).Value = "Fruta"
$ws.Cells.Item(14,7).Value = 100104
$ws.Cells.Item(14,8).Value = "Frutos de pepita"
$ws.Cells.Item(14,9).Value = 100104005
$ws.Cells.Item(14,10).Value = "Pera"
$ws.Cells.Item(14,11).Value = "Winter Nelis"
$ws.Cells.Item(14,12).Value = "Segunda"
$ws.Cells.Item(14,13).Value = 300
$ws.Cells.Item(14,14).Value = 17000
$ws.Cells.Item(14,15).Value = 18000
$ws.Cells.Item(14,16).Value = 17500
$ws.Cells.Item(14,17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(14,18).Value = "Región de O'Higgins"
$ws.Cells.Item(14,19).Value = 972
$ws.Cells.Item(14,20).Value = 18
